$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the statistical values in row 2 (regen sval data to filter save games)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.00000000002860089942657851
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 71518.04448319117
